$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 112183920
$ws.Range("Q2").Value = 763654
$ws.Range("R2").Value = 7448906
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()

# --- Row 3 ---
$ws.Range("Q3").Value = 763401
$ws.Range("R3").Value = 7448827
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()

# --- Row 4 ---
$ws.Range("A4").Value = 112183921
$ws.Range("B4").Value = 89405
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q4").Value = 763392
$ws.Range("R4").Value = 7448819
$ws.Range("Z4").ClearContents()
$ws.Range("AB4").ClearContents()

# --- Row 5 ---
$ws.Range("A5").Value = 112183036
$ws.Range("B5").Value = 89423
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 763482
$ws.Range("R5").Value = 7448939
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()

# --- Row 6 ---
$ws.Range("A6").Value = 112181983
$ws.Range("B6").Value = 89570
$ws.Range("E6").Value = 1588
$ws.Range("F6").Value = "Violmussling"
$ws.Range("G6").Value = "Trichaptum laricinum"
$ws.Range("H6").Value = "(P.Karst.) Ryvarden"
$ws.Range("Q6").Value = 763400
$ws.Range("R6").Value = 7448829
$ws.Range("Y6").Value = "'2023-07-06"
$ws.Range("Y6").ClearFormats()
$ws.Range("Z6").ClearContents()
$ws.Range("AA6").Value = "'2023-07-06"
$ws.Range("AA6").ClearFormats()
$ws.Range("AB6").ClearContents()

# --- Row 7 ---
$ws.Range("A7").Value = 112183947
$ws.Range("Q7").Value = 763391
$ws.Range("R7").Value = 7448820
$ws.Range("Y7").Value = "'2023-06-30"
$ws.Range("Y7").ClearFormats()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").Value = "'2023-06-30"
$ws.Range("AA7").ClearFormats()
$ws.Range("AB7").ClearContents()
